$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1417.0264
$ws.Range("I15").Value = 1417.0264
$ws.Range("K15").Value = 4251.0792
$ws.Range("M15").Value = -4082.0792
$ws.Range("H18").Value = 2499.375
$ws.Range("I18").Value = 2499.375
$ws.Range("K18").Value = 2499.375
$ws.Range("M18").Value = -2215.375
$ws.Range("H41").Value = 717.0833
$ws.Range("I41").Value = 720
$ws.Range("K41").Value = 720
$ws.Range("M41").Value = -280
$ws.Range("H75").Value = 272763740
$ws.Range("J75").Value = 272763740
$ws.Range("L75").Value = 272763740
$ws.Range("N75").Value = -272765612
$ws.Range("H78").Value = 272763740
$ws.Range("J78").Value = 272763740
$ws.Range("L78").Value = 818291220
$ws.Range("N78").Value = -818300580
$ws.Range("H92").Value = 1421.3158
$ws.Range("I92").Value = 1236.0714
$ws.Range("K92").Value = 1236.0714
$ws.Range("M92").Value = 11.92859999999996
$ws.Range("H116").Value = 12000.125
$ws.Range("I116").Value = 10199
$ws.Range("K116").Value = 10199
$ws.Range("M116").Value = -6757
$ws.Range("I135").Value = 4997
$ws.Range("K135").Value = 44973
$ws.Range("M135").Value = -42438
$ws.Range("H138").Value = 3181603.8
$ws.Range("J138").Value = 4836323
$ws.Range("L138").Value = 14508969
$ws.Range("N138").Value = -14519249
$ws.Range("H139").Value = 99803.86
$ws.Range("J139").Value = 99803.86
$ws.Range("L139").Value = 99803.86
$ws.Range("N139").Value = -110083.86
$ws.Range("H140").Value = 99999.25
$ws.Range("J140").Value = 99999.25
$ws.Range("L140").Value = 99999.25
$ws.Range("N140").Value = -110359.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2545.8333
$ws.Range("I2").Value = 2194.4614
$ws.Range("K2").Value = 2194.4614
$ws.Range("M2").Value = -2081.4614
$ws.Range("H4").Value = 315
$ws.Range("I4").Value = 278
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 278
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -162
$ws.Range("N4").Value = -732
$ws.Range("H5").Value = 169.82353
$ws.Range("I5").Value = 131.54546
$ws.Range("J5").Value = 240
$ws.Range("K5").Value = 131.54546
$ws.Range("L5").Value = 240
$ws.Range("M5").Value = -19.54545999999999
$ws.Range("N5").Value = -464
$ws.Range("H32").Value = 5215.9717
$ws.Range("I32").Value = 3658.4126
$ws.Range("K32").Value = 3658.4126
$ws.Range("M32").Value = -3371.4126
$ws.Range("H88").Value = 3162.4167
$ws.Range("J88").Value = 4500
$ws.Range("L88").Value = 4500
$ws.Range("N88").Value = -5312
$ws.Range("H91").Value = 3162.4167
$ws.Range("J91").Value = 4500
$ws.Range("L91").Value = 4500
$ws.Range("N91").Value = -7308
$ws.Range("H97").Value = 341.8
$ws.Range("I97").Value = 341.8
$ws.Range("K97").Value = 341.8
$ws.Range("M97").Value = 154.2
$ws.Range("H116").Value = 2545.8333
$ws.Range("I116").Value = 2194.4614
$ws.Range("K116").Value = 2194.4614
$ws.Range("M116").Value = 99.53859999999986
$ws.Range("H122").Value = 5292748.5
$ws.Range("I122").Value = 7408700
$ws.Range("J122").Value = 2869
$ws.Range("K122").Value = 22226100
$ws.Range("L122").Value = 8607
$ws.Range("M122").Value = -22223650
$ws.Range("N122").Value = -13507
$ws.Range("H131").Value = 72571.336
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 72571.336
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 72571.336
$ws.Range("N131").Value = -82651.336
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2545.8333
$ws.Range("I3").Value = 2194.4614
$ws.Range("K3").Value = 2194.4614
$ws.Range("M3").Value = -2080.4614
$ws.Range("H4").Value = 169.82353
$ws.Range("I4").Value = 131.54546
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 131.54546
$ws.Range("L4").Value = 240
$ws.Range("M4").Value = -16.54545999999999
$ws.Range("N4").Value = -470
$ws.Range("H86").Value = 29420876
$ws.Range("J86").Value = 500000000
$ws.Range("L86").Value = 500000000
$ws.Range("N86").Value = -500002246
$ws.Range("H89").Value = 29420876
$ws.Range("J89").Value = 500000000
$ws.Range("L89").Value = 2500000000
$ws.Range("N89").Value = -2500011232
$ws.Range("H94").Value = 2701.3333
$ws.Range("I94").Value = 2436.3333
$ws.Range("J94").Value = 2966.3333
$ws.Range("K94").Value = 2436.3333
$ws.Range("L94").Value = 2966.3333
$ws.Range("M94").Value = -1985.3333
$ws.Range("N94").Value = -3868.3333
$ws.Range("H99").Value = 3856.85
$ws.Range("I99").Value = 2846.3333
$ws.Range("J99").Value = 5372.625
$ws.Range("K99").Value = 2846.3333
$ws.Range("L99").Value = 5372.625
$ws.Range("M99").Value = -1348.3333
$ws.Range("N99").Value = -8368.625
$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 75000
$ws.Range("L108").Value = 75000
$ws.Range("N108").Value = -82680
$ws.Range("H134").Value = 2594.8057
$ws.Range("I134").Value = 2514.7856
$ws.Range("K134").Value = 7544.3568
$ws.Range("M134").Value = -5009.3568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 972.5
$ws.Range("I16").Value = 963.3333
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 963.3333
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -676.3333
$ws.Range("N16").Value = -1574
$ws.Range("H23").Value = 11740333
$ws.Range("I23").Value = 11740333
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 11740333
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -11740093
$ws.Range("H27").Value = 11740333
$ws.Range("I27").Value = 11740333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 11740333
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -11740141
$ws.Range("H105").Value = 3716.5
$ws.Range("I105").Value = 3656.7144
$ws.Range("J105").Value = 3800.2
$ws.Range("K105").Value = 3656.7144
$ws.Range("L105").Value = 3800.2
$ws.Range("M105").Value = -1909.7144
$ws.Range("N105").Value = -7294.2
$ws.Range("H107").Value = 850.15
$ws.Range("I107").Value = 914.8823
$ws.Range("J107").Value = 483.33334
$ws.Range("K107").Value = 914.8823
$ws.Range("L107").Value = 483.33334
$ws.Range("M107").Value = 1005.1177
$ws.Range("N107").Value = -4323.33334
$ws.Range("H113").Value = 972.5
$ws.Range("I113").Value = 963.3333
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 963.3333
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1206.6667
$ws.Range("N113").Value = -5340
$ws.Range("H137").Value = 98380
$ws.Range("J137").Value = 100975
$ws.Range("L137").Value = 100975
$ws.Range("N137").Value = -111175
$ws.Range("N23").ClearContents()
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.36842
$ws.Range("J2").Value = 211.57143
$ws.Range("L2").Value = 1269.42858
$ws.Range("N2").Value = -1495.42858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 33342876
$ws.Range("I113").Value = 71439300
$ws.Range("J113").Value = 8512.0625
$ws.Range("K113").Value = 71439300
$ws.Range("L113").Value = 8512.0625
$ws.Range("M113").Value = -71437130
$ws.Range("N113").Value = -12852.0625
$ws.Range("H123").Value = 64999
$ws.Range("J123").Value = 64999
$ws.Range("L123").Value = 64999
$ws.Range("N123").Value = -69899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2589.4666
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2589.4666
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2589.4666
$ws.Range("N46").Value = -2965.4666
$ws.Range("H55").Value = 1998.0588
$ws.Range("I55").Value = 1414.3334
$ws.Range("K55").Value = 1414.3334
$ws.Range("M55").Value = -1241.3334
$ws.Range("H61").Value = 4953.643
$ws.Range("I61").Value = 2370.6
$ws.Range("J61").Value = 6388.6665
$ws.Range("K61").Value = 2370.6
$ws.Range("L61").Value = 6388.6665
$ws.Range("M61").Value = -2168.6
$ws.Range("N61").Value = -6792.6665
$ws.Range("H113").Value = 4953.643
$ws.Range("I113").Value = 2370.6
$ws.Range("J113").Value = 6388.6665
$ws.Range("K113").Value = 2370.6
$ws.Range("L113").Value = 6388.6665
$ws.Range("M113").Value = -200.5999999999999
$ws.Range("N113").Value = -10728.6665
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 9999
$ws.Range("J18").Value = 9999
$ws.Range("L18").Value = 9999
$ws.Range("N18").Value = -10345
$ws.Range("H24").Value = 13999.667
$ws.Range("J24").Value = 13999.667
$ws.Range("L24").Value = 13999.667
$ws.Range("N24").Value = -14459.667
$ws.Range("H107").Value = 551.6539
$ws.Range("J107").Value = 607.2222
$ws.Range("L107").Value = 1821.6666
$ws.Range("N107").Value = -5661.6666
